$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value "D" below the existing B6:B8 data (A, B, C) into B9
$ws.Range("B9").Value = "D"

# Select the newly added cell, matching the saved selection state
$ws.Range("B9").Select()
